$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header row
$ws.Range("B1").Value = "Participant_ID"
$ws.Range("F1").Value = "Alignment score (Participant_Text|Participant_sim)"

# Source row
$ws.Range("B2").Value = "Text scenarios and sim scenarios (matched)"
$ws.Range("C2").Value = "Text scenarios"
$ws.Range("D2").Value = "Text scenarios"
$ws.Range("E2").Value = "Text scenarios"

# Definition row
$ws.Range("B3").Value = "Used to track and identify participants, also called Delegator ID"
$ws.Range("E3").Value = "Scenario presented to the participant in the text scenarios"
$ws.Range("F3").Value = "Comparse the KDMA measurement based on participant probe responses on the text scenario to participant probe responses in the simulated scenario"

# Selection + row height
$ws.Range("F4").Select()
$ws.Rows.Item(2).RowHeight = 45.75
